$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "5802202"
$ws.Range("A3").Value = "19499545"

$ws.Range("E9").Select()
